$wb = $excel.ActiveWorkbook

# Per-sheet last data row (header is row 1; data runs rows 2..lastRow)
$lastRows = @{ 1 = 7; 2 = 15; 3 = 16; 4 = 16 }
# Column E width (characters) to set explicitly for sheets that had one in the target
$colWidths = @{ 1 = 29.45; 2 = 22.93 }

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $lastRow = $lastRows[$i]

    # Header cell E1: same look as D1 ("link" header) -> copy its format, then set the text
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial(-4122)
    $ws.Range("E1").Value = "img"

    # Data cells E2:E<lastRow>: new plain style (based on D2's xf, font forced back to the
    # plain Arial 10 / black body font used across the sheet)
    $dataRange = "E2:E" + $lastRow
    $ws.Range("D2").Copy()
    $ws.Range($dataRange).PasteSpecial(-4122)
    $ws.Range($dataRange).Value = "./members/balidaan.jpg"
    $ws.Range($dataRange).Font.Name = "Arial"
    $ws.Range($dataRange).Font.Size = 10
    $ws.Range($dataRange).Font.Color = 0

    if ($colWidths.ContainsKey($i)) {
        $ws.Columns.Item(5).ColumnWidth = $colWidths[$i]
    }
}

$excel.CutCopyMode = 0

# Restore each sheet's own cursor position, finishing on Sheet1 so it ends up the active tab
$wb.Worksheets.Item(2).Range("C36").Select()
$wb.Worksheets.Item(3).Range("C34").Select()
$wb.Worksheets.Item(4).Range("D28").Select()
$wb.Worksheets.Item(1).Range("E10").Select()
